$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "XGB"
$ws.Range("C2").Value = "Reduced"
$ws.Range("D2").Value = "Default"
$ws.Range("E2").Value = 0.9468463503832474
$ws.Range("F2").Value = 0.9423942394239424
$ws.Range("G2").Value = 0.9510496568429552
$ws.Range("H2").Value = 0.9910657197018637
$ws.Range("I2").Value = 9423
$ws.Range("J2").Value = 9477
$ws.Range("K2").Value = 485
$ws.Range("L2").Value = 576
$ws.Range("M2").Value = 292850
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Random Forest"
$ws.Range("C3").Value = "Reduced"
$ws.Range("D3").Value = "Default"
$ws.Range("E3").Value = 0.9469966434547368
$ws.Range("F3").Value = 0.9681968196819682
$ws.Range("G3").Value = 0.9289895403512139
$ws.Range("H3").Value = 0.9908716278172688
$ws.Range("I3").Value = 9681
$ws.Range("J3").Value = 9222
$ws.Range("K3").Value = 740
$ws.Range("L3").Value = 318
$ws.Range("M3").Value = 166400
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Decision Tree"
$ws.Range("C4").Value = "Reduced"
$ws.Range("D4").Value = "Best"
$ws.Range("E4").Value = 0.8570712890135764
$ws.Range("F4").Value = 0.8905890589058906
$ws.Range("G4").Value = 0.8350525131282821
$ws.Range("H4").Value = 0.9347231601297051
$ws.Range("I4").Value = 8905
$ws.Range("J4").Value = 8203
$ws.Range("K4").Value = 1759
$ws.Range("L4").Value = 1094
$ws.Range("M4").Value = 564590
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Logistic Regression"
$ws.Range("C5").Value = "Full"
$ws.Range("D5").Value = "Default"
$ws.Range("E5").Value = 0.8232553479284604
$ws.Range("F5").Value = 0.8277827782778278
$ws.Range("G5").Value = 0.8208866408806903
$ws.Range("H5").Value = 0.8929961456294193
$ws.Range("I5").Value = 8277
$ws.Range("J5").Value = 8156
$ws.Range("K5").Value = 1806
$ws.Range("L5").Value = 1722
$ws.Range("M5").Value = 879060
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Gaussian Naive-Bayes"
$ws.Range("C6").Value = "Full"
$ws.Range("D6").Value = "Default"
$ws.Range("E6").Value = 0.7034717699514053
$ws.Range("F6").Value = 0.8596859685968596
$ws.Range("G6").Value = 0.6555826723611958
$ws.Range("H6").Value = 0.8058707296146197
$ws.Range("I6").Value = 8596
$ws.Range("J6").Value = 5446
$ws.Range("K6").Value = 4516
$ws.Range("L6").Value = 1403
$ws.Range("M6").Value = 746660
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "Bernoulli Naive-Bayes"
$ws.Range("C7").Value = "Full"
$ws.Range("D7").Value = "Default"
$ws.Range("E7").Value = 0.6491658734532338
$ws.Range("F7").Value = 0.6465646564656465
$ws.Range("G7").Value = 0.6507952486410308
$ws.Range("H7").Value = 0.7062773432532974
$ws.Range("I7").Value = 6465
$ws.Range("J7").Value = 6493
$ws.Range("K7").Value = 3469
$ws.Range("L7").Value = 3534
$ws.Range("M7").Value = 1801690
